# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" on every sheet.
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   are refreshed to reflect the new handoff report run.
# - Because "Ready for handoff" is wider than "In Translation", the Status
#   column on each sheet is widened to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# Width (in COM ColumnWidth units) that reproduces the widened Status column.
$statusColWidth = 16.333333333333336

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-24 22:59:02"
$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-24 22:58:56"
$zhcn.Columns.Item(3).ColumnWidth = $statusColWidth

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-24 22:59:02"
$dede.Columns.Item(3).ColumnWidth = $statusColWidth
